$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row containing StockCode 71053 (original row 3). This shifts
# rows 4-7 up by one, so the "missing value" sample becomes cleaner.
$ws.Rows.Item(3).Delete()

# In the resulting row 4 (StockCode 84029G), clear out the CustomerID value
# to create the sample with a missing value, as described in the commit
# message.
$ws.Range("F4").ClearContents()

# Update the used range dimension / selection to match a plain view.
$ws.Range("A1").Select()
